$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.046.08'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.51%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.645.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.32%  '

$ws.Range('E4').Value = '  -0.15%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.14%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5103'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.29%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2564'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.25%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06360'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.48%  '

$ws.Range('E10').Value = '  -0.20%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07766'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.17%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.295'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.41%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.650.27'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.22%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5442'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.16%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.32'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.54%  '

$ws.Range('E16').Value = '  -1.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.074.71'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.002'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.26%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '199.02'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.432'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.01%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.941'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.07%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.050'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.22%  '

$ws.Range('E23').Value = '  -0.21%  '

$ws.Range('E24').Value = '  -0.24%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.06%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1194'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.40%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.816'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.62%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.64'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.44%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.236'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.53%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04861'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.70%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.260'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '

$ws.Range('E32').Value = '  -0.79%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.526'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.30%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.368'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9000'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.82%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.586'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.66%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.142.91'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.32%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5465'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.54%  '

$ws.Range('E39').Value = '  +0.29%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.002'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.32%  '

$ws.Range('B41').Value = 'BabyDogeCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₈130'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +8.97%  '

$ws.Range('B42').Value = 'mCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.529'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.24%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8127'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.52%  '

$ws.Range('E44').Value = '  -0.22%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.393'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.23%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.782.95'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.37%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4530'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.24%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.02'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.51%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.000'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.90%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05057'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.55%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.003'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.41%  '
